# Update progress ("진행상태") percentages in the 개발목록 (dev task list) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11, 13-16: DB 연결 task block progress 60% -> 90%
$ws.Range("I11").Value = 0.9
$ws.Range("I13").Value = 0.9
$ws.Range("I14").Value = 0.9
$ws.Range("I15").Value = 0.9
$ws.Range("I16").Value = 0.9

# Row 26: 관리자 섹션 progress 30% -> 90%
$ws.Range("I26").Value = 0.9

# Row 28: 회원관리 페이지 progress 60% -> 90%
$ws.Range("I28").Value = 0.9

# Row 29: progress 10% -> 20%
$ws.Range("I29").Value = 0.2

# Move the active selection down to I30, matching the saved cursor position.
$ws.Range("I30").Select()
